# Pokedex.xlsx update:
#  - Insert a new "Sheet1" (type-effectiveness chart) after "Branch Evos"
#  - National sheet: rename "Method(s)" header to "Means", drop a stray
#    G37 value, and refocus the view
#  - Branch Evos sheet: swap the Branches/Level columns, rework the
#    Evolution/Method columns, and refocus the view

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. Insert a new "Sheet1" (type effectiveness chart) after "Branch Evos"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("Branch Evos")
$typesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)

$typesSheet.Range("A1").Value = "Attacking"
$typesSheet.Range("B1").Value = "Normal"
$typesSheet.Range("C1").Value = "Fire"
$typesSheet.Range("D1").Value = "Water"
$typesSheet.Range("E1").Value = "Electric"
$typesSheet.Range("F1").Value = "Grass"
$typesSheet.Range("G1").Value = "Ice"
$typesSheet.Range("H1").Value = "Fighting"
$typesSheet.Range("I1").Value = "Poison"
$typesSheet.Range("J1").Value = "Ground"
$typesSheet.Range("K1").Value = "Flying"
$typesSheet.Range("L1").Value = "Psychic"
$typesSheet.Range("M1").Value = "Bug"
$typesSheet.Range("N1").Value = "Rock"
$typesSheet.Range("O1").Value = "Ghost"
$typesSheet.Range("P1").Value = "Dragon"
$typesSheet.Range("Q1").Value = "Dark"
$typesSheet.Range("R1").Value = "Steel"
$typesSheet.Range("S1").Value = "Fairy"
$typesSheet.Range("A1:S1").Font.Bold = $true

$typesSheet.Range("A3").Value = "Normal"
$typesSheet.Range("A3").Font.Bold = $true
$typesSheet.Range("N3").Value = 0.5
$typesSheet.Range("O3").Value = 0
$typesSheet.Range("R3").Value = 0.5

$typesSheet.Range("A4").Value = "Fire"
$typesSheet.Range("A4").Font.Bold = $true
$typesSheet.Range("C4").Value = 0.5
$typesSheet.Range("D4").Value = 0.5
$typesSheet.Range("F4").Value = 2
$typesSheet.Range("G4").Value = 2
$typesSheet.Range("M4").Value = 2
$typesSheet.Range("N4").Value = 0.5
$typesSheet.Range("P4").Value = 0.5
$typesSheet.Range("R4").Value = 2

$typesSheet.Range("A5").Value = "Water"
$typesSheet.Range("A5").Font.Bold = $true
$typesSheet.Range("C5").Value = 2
$typesSheet.Range("D5").Value = 0.5
$typesSheet.Range("F5").Value = 0.5
$typesSheet.Range("J5").Value = 2
$typesSheet.Range("N5").Value = 2
$typesSheet.Range("P5").Value = 0.5

$typesSheet.Range("A6").Value = "Electric"
$typesSheet.Range("A6").Font.Bold = $true
$typesSheet.Range("D6").Value = 2
$typesSheet.Range("E6").Value = 0.5
$typesSheet.Range("F6").Value = 0.5
$typesSheet.Range("J6").Value = 0
$typesSheet.Range("K6").Value = 2
$typesSheet.Range("P6").Value = 0.5

$typesSheet.Range("A7").Value = "Grass"
$typesSheet.Range("A7").Font.Bold = $true
$typesSheet.Range("C7").Value = 0.5
$typesSheet.Range("D7").Value = 2
$typesSheet.Range("F7").Value = 0.5
$typesSheet.Range("I7").Value = 0.5
$typesSheet.Range("J7").Value = 2
$typesSheet.Range("K7").Value = 0.5
$typesSheet.Range("M7").Value = 0.5
$typesSheet.Range("N7").Value = 2
$typesSheet.Range("P7").Value = 0.5
$typesSheet.Range("R7").Value = 0.5

$typesSheet.Range("A8").Value = "Ice"
$typesSheet.Range("A8").Font.Bold = $true
$typesSheet.Range("C8").Value = 0.5
$typesSheet.Range("D8").Value = 0.5
$typesSheet.Range("F8").Value = 2
$typesSheet.Range("G8").Value = 0.5
$typesSheet.Range("J8").Value = 2
$typesSheet.Range("K8").Value = 2
$typesSheet.Range("P8").Value = 2
$typesSheet.Range("R8").Value = 0.5

$typesSheet.Range("A9").Value = "Fighting"
$typesSheet.Range("A9").Font.Bold = $true
$typesSheet.Range("B9").Value = 2
$typesSheet.Range("G9").Value = 2
$typesSheet.Range("I9").Value = 0.5
$typesSheet.Range("K9").Value = 0.5
$typesSheet.Range("L9").Value = 0.5
$typesSheet.Range("M9").Value = 0.5
$typesSheet.Range("N9").Value = 2
$typesSheet.Range("O9").Value = 0
$typesSheet.Range("Q9").Value = 2
$typesSheet.Range("R9").Value = 2
$typesSheet.Range("S9").Value = 0.5

$typesSheet.Range("A10").Value = "Poison"
$typesSheet.Range("A10").Font.Bold = $true
$typesSheet.Range("F10").Value = 2
$typesSheet.Range("I10").Value = 0.5
$typesSheet.Range("J10").Value = 0.5
$typesSheet.Range("N10").Value = 0.5
$typesSheet.Range("O10").Value = 0.5
$typesSheet.Range("R10").Value = 0
$typesSheet.Range("S10").Value = 2

$typesSheet.Range("A11").Value = "Ground"
$typesSheet.Range("A11").Font.Bold = $true
$typesSheet.Range("C11").Value = 2
$typesSheet.Range("E11").Value = 2
$typesSheet.Range("F11").Value = 0.5
$typesSheet.Range("I11").Value = 2
$typesSheet.Range("K11").Value = 0
$typesSheet.Range("M11").Value = 0.5
$typesSheet.Range("N11").Value = 2
$typesSheet.Range("R11").Value = 2

$typesSheet.Range("A12").Value = "Flying"
$typesSheet.Range("A12").Font.Bold = $true
$typesSheet.Range("E12").Value = 0.5
$typesSheet.Range("F12").Value = 2
$typesSheet.Range("H12").Value = 2
$typesSheet.Range("M12").Value = 2
$typesSheet.Range("N12").Value = 0.5
$typesSheet.Range("R12").Value = 0.5

$typesSheet.Range("A13").Value = "Psychic"
$typesSheet.Range("A13").Font.Bold = $true
$typesSheet.Range("H13").Value = 2
$typesSheet.Range("I13").Value = 2
$typesSheet.Range("L13").Value = 0.5
$typesSheet.Range("Q13").Value = 0
$typesSheet.Range("R13").Value = 0.5

$typesSheet.Range("A14").Value = "Bug"
$typesSheet.Range("A14").Font.Bold = $true
$typesSheet.Range("C14").Value = 0.5
$typesSheet.Range("F14").Value = 2
$typesSheet.Range("H14").Value = 0.5
$typesSheet.Range("I14").Value = 0.5
$typesSheet.Range("K14").Value = 0.5
$typesSheet.Range("L14").Value = 2
$typesSheet.Range("O14").Value = 0.5
$typesSheet.Range("Q14").Value = 2
$typesSheet.Range("R14").Value = 0.5
$typesSheet.Range("S14").Value = 0.5

$typesSheet.Range("A15").Value = "Rock"
$typesSheet.Range("A15").Font.Bold = $true
$typesSheet.Range("C15").Value = 2
$typesSheet.Range("G15").Value = 2
$typesSheet.Range("H15").Value = 0.5
$typesSheet.Range("J15").Value = 0.5
$typesSheet.Range("K15").Value = 2
$typesSheet.Range("M15").Value = 2
$typesSheet.Range("R15").Value = 0.5

$typesSheet.Range("A16").Value = "Ghost"
$typesSheet.Range("A16").Font.Bold = $true
$typesSheet.Range("B16").Value = 0
$typesSheet.Range("L16").Value = 2
$typesSheet.Range("O16").Value = 2
$typesSheet.Range("Q16").Value = 0.5

$typesSheet.Range("A17").Value = "Dragon"
$typesSheet.Range("A17").Font.Bold = $true
$typesSheet.Range("P17").Value = 2
$typesSheet.Range("R17").Value = 0.5
$typesSheet.Range("S17").Value = 0

$typesSheet.Range("A18").Value = "Dark"
$typesSheet.Range("A18").Font.Bold = $true
$typesSheet.Range("H18").Value = 0.5
$typesSheet.Range("L18").Value = 2
$typesSheet.Range("O18").Value = 2
$typesSheet.Range("Q18").Value = 0.5
$typesSheet.Range("S18").Value = 0.5

$typesSheet.Range("A19").Value = "Steel"
$typesSheet.Range("A19").Font.Bold = $true
$typesSheet.Range("C19").Value = 0.5
$typesSheet.Range("D19").Value = 0.5
$typesSheet.Range("E19").Value = 0.5
$typesSheet.Range("G19").Value = 2
$typesSheet.Range("N19").Value = 2
$typesSheet.Range("R19").Value = 0.5
$typesSheet.Range("S19").Value = 2

$typesSheet.Range("A20").Value = "Fairy"
$typesSheet.Range("A20").Font.Bold = $true
$typesSheet.Range("C20").Value = 0.5
$typesSheet.Range("H20").Value = 2
$typesSheet.Range("I20").Value = 0.5
$typesSheet.Range("P20").Value = 2
$typesSheet.Range("Q20").Value = 2
$typesSheet.Range("R20").Value = 0.5

# row 2 is intentionally blank except for a styled, empty A2
$typesSheet.Range("A2").Font.Bold = $true

# ---------------------------------------------------------------------
# 1. National (sheet 1): header + stray-cell cleanup
# ---------------------------------------------------------------------
$national = $wb.Worksheets.Item("National")
$national.Range("H1").Value = "Means"
$national.Range("G37").ClearContents()

# ---------------------------------------------------------------------
# 2. Branch Evos (sheet 2): column rework
# ---------------------------------------------------------------------
$branch = $wb.Worksheets.Item("Branch Evos")
$branch.Range("A1:E15").ClearContents()

$branch.Range("A1").Value = "Name"
$branch.Range("B1").Value = "Level"
$branch.Range("C1").Value = "Branches"
$branch.Range("D1").Value = "Evolution"
$branch.Range("E1").Value = "Method"
$branch.Range("A1:E1").Font.Bold = $true

$branch.Range("A2").Value = "Gloom"
$branch.Range("C2").Value = 2
$branch.Range("D2").Value = "Vileplume"
$branch.Range("E2").Value = "Leaf Stone"

$branch.Range("D3").Value = "Bellossom"
$branch.Range("E3").Value = "Sun Stone"

$branch.Range("A4").Value = "Poliwhirl"
$branch.Range("C4").Value = 2
$branch.Range("D4").Value = "Poliwrath"
$branch.Range("E4").Value = "Water Stone"

$branch.Range("D5").Value = "Politoed"
$branch.Range("E5").Value = "Trade with King's Rock"

$branch.Range("A6").Value = "Slowpoke"
$branch.Range("C6").Value = 2
$branch.Range("D6").Value = "Slowbro"
$branch.Range("E6").Value = 37

$branch.Range("D7").Value = "Slowking"
$branch.Range("E7").Value = "Trade with King's Rock"

$branch.Range("A8").Value = "Eevee"
$branch.Range("C8").Value = 5
$branch.Range("D8").Value = "Vaporeon"
$branch.Range("E8").Value = "Water Stone"

$branch.Range("D9").Value = "Jolteon"
$branch.Range("E9").Value = "Thunder Stone"

$branch.Range("D10").Value = "Flareon"
$branch.Range("E10").Value = "Fire Stone"

$branch.Range("D11").Value = "Espeon"
$branch.Range("E11").Value = "Happiness during Day"

$branch.Range("D12").Value = "Umbreon"
$branch.Range("E12").Value = "Happiness during Night"

$branch.Range("A13").Value = "Tyrogue"
$branch.Range("B13").Value = 20
$branch.Range("C13").Value = 3
$branch.Range("D13").Value = "Hitmonlee"
$branch.Range("E13").Value = "Attack > Defense"

$branch.Range("D14").Value = "Hitmonchan"
$branch.Range("E14").Value = "Attack < Defense"

$branch.Range("D15").Value = "Hitmontop"
$branch.Range("E15").Value = "Attack = Defense"
